# ECU1_EntranceGate_TestCases.xlsx - add two new test cases (TC_Gate_16 / TC_Gate_17)
# for "Validate functionality of LEDs (Entrnce Gate)" scenario, following the same
# layout pattern used for the existing scenario blocks (rows 20-21 for LCD/Buzzer).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Build rows 22 (blank separator), 23 and 24 (new test case rows) by copying
#    the formatting of the analogous existing rows (19 / 20 / 21) and then
#    filling in the new values.
# ---------------------------------------------------------------------------

$ws.Range("A19:M19").Copy() | Out-Null
$ws.Range("A22:M22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Both new content rows start from row 21's formatting (no top border on A,
# wrapped G/H) ...
$ws.Range("A21:M21").Copy() | Out-Null
$ws.Range("A23:M23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A24:M24").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ... then row 23's A cell (first row of the new merged scenario block) gets
# the bordered "group header" look, matching A20/A15/A10/A8.
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$ws.Application.CutCopyMode = 0

# Row heights (thickBot separator + the two content rows)
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 120.75
$ws.Rows.Item(24).RowHeight = 120

# ---------------------------------------------------------------------------
# 2) Values for the two new test cases - entered column-by-column (both rows
#    of a column filled before moving to the next column) which is the order
#    the shared-strings table needs to come out in.
# ---------------------------------------------------------------------------
$ws.Range("A23").Value2 = "Validate functionality of LEDs  (Entrnce Gate)"

$ws.Range("B23").Value2 = "TC_Gate_16"
$ws.Range("B24").Value2 = "TC_Gate_17"

$ws.Range("C23").Value2 = "Validate that Green Led "
$ws.Range("C24").Value2 = "Validate that Red Led "

$ws.Range("E23").Value2 = "1-Atmel Studio `n2- Proteus Simulation`n3- RFID Driver `n4- UART Driver`n5- SPI Driver`n6- LCD Driver `n7- Buzzer`n"
$ws.Range("E24").Value2 = "1-Atmel Studio `n2- Proteus Simulation`n3- RFID Driver `n4- UART Driver`n5- SPI Driver`n6- LCD Driver `n7- Buzzer`n"

$ws.Range("F23").Value2 = "Enter Valid Driver`n1- username = ""Mohamed""`n2- ID= ""0000001""`n"
$ws.Range("F24").Value2 = "Enter invalid Driver`n1- username = ""Md""`n2- ID= ""0000001""`n"

$ws.Range("G23").Value2 = "Green Led on `nRed Led off"
$ws.Range("H23").Value2 = "Green Led on `nRed Led off"
$ws.Range("G24").Value2 = "Green Led off `nRed Led on"
$ws.Range("H24").Value2 = "Green Led off `nRed Led on"

$ws.Range("I23").Value2 = "Pass"
$ws.Range("J23").Value2 = "Mohamed Abd El-Naby"
$ws.Range("K23").Value2 = "Mohamed Abd El-Naby"
$ws.Range("L23").Value2 = "Functional test"

$ws.Range("I24").Value2 = "Pass"
$ws.Range("J24").Value2 = "Mohamed Abd El-Naby"
$ws.Range("K24").Value2 = "Mohamed Abd El-Naby"
$ws.Range("L24").Value2 = "Functional test"

# ---------------------------------------------------------------------------
# 4) Merge the scenario-objective cell across the two new rows
# ---------------------------------------------------------------------------
$ws.Range("A23:A24").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 5) Conditional formatting: reproduce the same 4-rule (inProgress / onGoing /
#    Pass / Fail) "Highlight Cell" status formatting used on every other
#    scenario block, for the sqref groups belonging to the new rows.
# ---------------------------------------------------------------------------
function Add-StatusConditionalFormat($range) {
    $fc = $range.FormatConditions

    $c1 = $fc.Add(1, 3, '="inProgress"')
    $c1.Font.Color = 393372
    $c1.Interior.Color = 13551615

    $c2 = $fc.Add(1, 3, '="onGoing"')
    $c2.Font.Color = 24832
    $c2.Interior.Color = 13561798

    $c3 = $fc.Add(1, 3, '="Pass"')
    $c3.Font.Color = 22428
    $c3.Interior.Color = 10284031

    $c4 = $fc.Add(1, 3, '="Fail"')
    $c4.Font.Color = 22428
    $c4.Interior.Color = 10284031
}

Add-StatusConditionalFormat($ws.Range("A22:M22"))
Add-StatusConditionalFormat($ws.Range("A23"))
Add-StatusConditionalFormat($ws.Range("B23:B24"))
Add-StatusConditionalFormat($ws.Range("E23:E24"))
Add-StatusConditionalFormat($ws.Range("F23:F24"))
Add-StatusConditionalFormat($ws.Range("I23"))
Add-StatusConditionalFormat($ws.Range("I24"))
Add-StatusConditionalFormat($ws.Range("J23:K24"))

# ---------------------------------------------------------------------------
# 6) Sheet view - scroll position / selection moved to reflect the new rows
# ---------------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A20"), $true)
$ws.Range("C26").Select() | Out-Null

Write-Host "Done adding TC_Gate_16 / TC_Gate_17"
